# Apply the StructureDefinition metadata update:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date updated
#  - Publisher contact info replaced with Publisher name + Jurisdiction
#  - Duplicate "Contact" row removed
#  - Elements sheet root extension Short/Definition filled in with the
#    extension's actual name/description instead of generic placeholders

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value2 = "6.0.0"

# Date: refreshed publication timestamp
$wsMeta.Range("B8").Value2 = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$wsMeta.Range("B9").Value2 = "Alvearie Team"

# Former duplicate "Contact" row repurposed into "Jurisdiction"
$wsMeta.Range("A10").Value2 = "Jurisdiction"
$wsMeta.Range("B10").Value2 = "United States of America"

# Remove the now-redundant second "Contact" / "No display for ContactDetail" row
$wsMeta.Rows.Item(11).Delete()

$wsElements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition updated to the extension-specific text
$wsElements.Range("K2").Value2 = "CareGapComplianceMet"
$wsElements.Range("L2").Value2 = "Indicates if the patient received the targeted treatment to close the care gap."
